$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 573.375
$ws.Range("J41").Value = 550
$ws.Range("L41").Value = 550
$ws.Range("N41").Value = -1430

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()

$ws.Range("H107").Value = 327.4
$ws.Range("I107").Value = 359.5
$ws.Range("K107").Value = 359.5
$ws.Range("M107").Value = 1560.5

$ws.Range("H112").Value = 1647.5714
$ws.Range("J112").Value = 1794.3889
$ws.Range("L112").Value = 5383.1667
$ws.Range("N112").Value = -7599.1667

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H137").Value = 1673.8
$ws.Range("I137").Value = 1205.3334
$ws.Range("K137").Value = 3616.0002
$ws.Range("M137").Value = -1066.0002

$ws.Range("H138").Value = 2500.4736
$ws.Range("J138").Value = 2771.889
$ws.Range("L138").Value = 8315.667000000001
$ws.Range("N138").Value = -18595.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1234.2858
$ws.Range("I2").Value = 1327
$ws.Range("K2").Value = 1327
$ws.Range("M2").Value = -1214

$ws.Range("H45").Value = 2169.111
$ws.Range("I45").Value = 2190.25
$ws.Range("K45").Value = 2190.25
$ws.Range("M45").Value = -1813.25

$ws.Range("H97").Value = 350.33334
$ws.Range("I97").Value = 227
$ws.Range("J97").Value = 1337
$ws.Range("K97").Value = 227
$ws.Range("L97").Value = 1337
$ws.Range("M97").Value = 269
$ws.Range("N97").Value = -2329

$ws.Range("H102").Value = 1286
$ws.Range("I102").Value = 1415.6666
$ws.Range("J102").Value = 119
$ws.Range("K102").Value = 1415.6666
$ws.Range("L102").Value = 119
$ws.Range("M102").Value = 206.3334
$ws.Range("N102").Value = -3363

$ws.Range("H116").Value = 1234.2858
$ws.Range("I116").Value = 1327
$ws.Range("K116").Value = 1327
$ws.Range("M116").Value = 967

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H139").Value = 74749.75
$ws.Range("J139").Value = 74749.75
$ws.Range("L139").Value = 74749.75
$ws.Range("N139").Value = -85029.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1234.2858
$ws.Range("I3").Value = 1327
$ws.Range("K3").Value = 1327
$ws.Range("M3").Value = -1213

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 50000
$ws.Range("J53").Value = 50000
$ws.Range("L53").Value = 50000
$ws.Range("N53").Value = -51214

$ws.Range("H68").Value = 69996.664
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 69996.664
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 69996.664
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -71494.664

$ws.Range("H71").Value = 69996.664
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 69996.664
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 209989.992
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -217477.992

$ws.Range("H134").Value = 3651.5386
$ws.Range("I134").Value = 3664.4443
$ws.Range("J134").Value = 3622.5
$ws.Range("K134").Value = 10993.3329
$ws.Range("L134").Value = 10867.5
$ws.Range("M134").Value = -8458.332900000001
$ws.Range("N134").Value = -15937.5

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 53657.43
$ws.Range("J141").Value = 53657.43
$ws.Range("L141").Value = 53657.43
$ws.Range("N141").Value = -64017.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 6173.25
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H56").Value = 18418.559
$ws.Range("I56").Value = 18418.559
$ws.Range("K56").Value = 18418.559
$ws.Range("M56").Value = -17888.559

$ws.Range("H75").Value = 923.5
$ws.Range("J75").Value = 1347
$ws.Range("L75").Value = 4041
$ws.Range("N75").Value = -6037

$ws.Range("H78").Value = 923.5
$ws.Range("J78").Value = 1347
$ws.Range("L78").Value = 12123
$ws.Range("N78").Value = -22107

$ws.Range("H87").Value = 2831.6667
$ws.Range("I87").Value = 2831.6667
$ws.Range("K87").Value = 8495.000100000001
$ws.Range("M87").Value = -7247.000100000001

$ws.Range("H88").Value = 2500
$ws.Range("I88").Value = 2500
$ws.Range("K88").Value = 7500
$ws.Range("M88").Value = -7072

$ws.Range("H90").Value = 2831.6667
$ws.Range("I90").Value = 2831.6667
$ws.Range("K90").Value = 25485.0003
$ws.Range("M90").Value = -19245.0003

$ws.Range("H91").Value = 2500
$ws.Range("I91").Value = 2500
$ws.Range("K91").Value = 7500
$ws.Range("M91").Value = -6018

$ws.Range("H97").Value = 1500
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 1404.625
$ws.Range("I107").Value = 5050
$ws.Range("J107").Value = 189.5
$ws.Range("K107").Value = 15150
$ws.Range("L107").Value = 568.5
$ws.Range("M107").Value = -13230
$ws.Range("N107").Value = -4408.5

$ws.Range("H117").Value = 1607.8334
$ws.Range("J117").Value = 1966
$ws.Range("L117").Value = 5898
$ws.Range("N117").Value = -12782

$ws.Range("H121").Value = 2285.4
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2285.4
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 6856.200000000001
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -9476.200000000001

$ws.Range("H122").Value = 78530.84
$ws.Range("J122").Value = 144560.86
$ws.Range("L122").Value = 1301047.74
$ws.Range("N122").Value = -1305947.74

$ws.Range("H132").Value = 4531
$ws.Range("I132").Value = 5500
$ws.Range("J132").Value = 4208
$ws.Range("K132").Value = 49500
$ws.Range("L132").Value = 37872
$ws.Range("M132").Value = -46970
$ws.Range("N132").Value = -42932

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2660.9
$ws.Range("I132").Value = 1885.5454
$ws.Range("K132").Value = 5656.6362
$ws.Range("M132").Value = -3126.6362

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2000
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1639
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 2000
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -752
$ws.Range("N85").ClearContents()

$ws.Range("H132").Value = 4343.222
$ws.Range("I132").Value = 3964.3333
$ws.Range("J132").Value = 4532.6665
$ws.Range("K132").Value = 11892.9999
$ws.Range("L132").Value = 13597.9995
$ws.Range("M132").Value = -9362.999899999999
$ws.Range("N132").Value = -18657.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1269.6666
$ws.Range("I100").Value = 1269.6666
$ws.Range("K100").Value = 2539.3332
$ws.Range("M100").Value = -1998.3332
